$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 821599
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 821599
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 821599
$ws.Range("N43").Value = -821737

$ws.Range("H86").Value = 92595464
$ws.Range("I86").Value = 111113240
$ws.Range("J86").Value = 18524350
$ws.Range("K86").Value = 111113240
$ws.Range("L86").Value = 18524350
$ws.Range("M86").Value = -111112117
$ws.Range("N86").Value = -18526596

$ws.Range("H89").Value = 92595464
$ws.Range("I89").Value = 111113240
$ws.Range("J89").Value = 18524350
$ws.Range("K89").Value = 555566200
$ws.Range("L89").Value = 92621750
$ws.Range("M89").Value = -555560584
$ws.Range("N89").Value = -92632982

$ws.Range("H132").Value = 1591.7291
$ws.Range("I132").Value = 1195.6316
$ws.Range("J132").Value = 3096.9
$ws.Range("K132").Value = 3586.8948
$ws.Range("L132").Value = 9290.700000000001
$ws.Range("M132").Value = -1056.8948
$ws.Range("N132").Value = -14350.7

$ws.Range("H138").Value = 6275.6333
$ws.Range("I138").Value = 1084.8667
$ws.Range("J138").Value = 11466.4
$ws.Range("K138").Value = 3254.6001
$ws.Range("L138").Value = 34399.2
$ws.Range("M138").Value = 1885.3999
$ws.Range("N138").Value = -44679.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1868204.5
$ws.Range("I32").Value = 1986530.1
$ws.Range("J32").Value = 4575.25
$ws.Range("K32").Value = 1986530.1
$ws.Range("L32").Value = 4575.25
$ws.Range("M32").Value = -1986243.1
$ws.Range("N32").Value = -5149.25

$ws.Range("H63").Value = 3499.4
$ws.Range("I63").Value = 4248.5
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 4248.5
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -3562.5
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 3499.4
$ws.Range("I66").Value = 4248.5
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 21242.5
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -17810.5
$ws.Range("N66").Value = -21864

$ws.Range("H102").Value = 1261.5385
$ws.Range("I102").Value = 1225
$ws.Range("J102").Value = 1462.5
$ws.Range("K102").Value = 1225
$ws.Range("L102").Value = 1462.5
$ws.Range("M102").Value = 397
$ws.Range("N102").Value = -4706.5

$ws.Range("H122").Value = 12266.833
$ws.Range("I122").Value = 16838.846
$ws.Range("J122").Value = 6863.5454
$ws.Range("K122").Value = 50516.538
$ws.Range("L122").Value = 20590.6362
$ws.Range("M122").Value = -48066.538
$ws.Range("N122").Value = -25490.6362

$ws.Range("H132").Value = 7283.3335
$ws.Range("I132").Value = 6002.2915
$ws.Range("J132").Value = 9845.416999999999
$ws.Range("K132").Value = 18006.8745
$ws.Range("L132").Value = 29536.251
$ws.Range("M132").Value = -15476.8745
$ws.Range("N132").Value = -34596.251

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 76796.60000000001
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 76796.60000000001
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 76796.60000000001
$ws.Range("N81").Value = -78918.60000000001

$ws.Range("H84").Value = 76796.60000000001
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 76796.60000000001
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 230389.8
$ws.Range("N84").Value = -240997.8

$ws.Range("H99").Value = 2599587.5
$ws.Range("I99").Value = 1949.05
$ws.Range("J99").Value = 6063105.5
$ws.Range("K99").Value = 1949.05
$ws.Range("L99").Value = 6063105.5
$ws.Range("M99").Value = -451.05
$ws.Range("N99").Value = -6066101.5

$ws.Range("H105").Value = 1812.5
$ws.Range("I105").Value = 1261
$ws.Range("J105").Value = 2854.2222
$ws.Range("K105").Value = 1261
$ws.Range("L105").Value = 2854.2222
$ws.Range("M105").Value = 486
$ws.Range("N105").Value = -6348.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4733.1665
$ws.Range("I16").Value = 3091.6667
$ws.Range("J16").Value = 6374.6665
$ws.Range("K16").Value = 3091.6667
$ws.Range("L16").Value = 6374.6665
$ws.Range("M16").Value = -2804.6667
$ws.Range("N16").Value = -6948.6665

$ws.Range("H31").Value = 7536
$ws.Range("I31").Value = 2331.5715
$ws.Range("J31").Value = 11739.577
$ws.Range("K31").Value = 2331.5715
$ws.Range("L31").Value = 11739.577
$ws.Range("M31").Value = -2036.5715
$ws.Range("N31").Value = -12329.577

$ws.Range("H34").Value = 7536
$ws.Range("I34").Value = 2331.5715
$ws.Range("J34").Value = 11739.577
$ws.Range("K34").Value = 2331.5715
$ws.Range("L34").Value = 11739.577
$ws.Range("M34").Value = -2129.5715
$ws.Range("N34").Value = -12143.577

$ws.Range("H113").Value = 4733.1665
$ws.Range("I113").Value = 3091.6667
$ws.Range("J113").Value = 6374.6665
$ws.Range("K113").Value = 3091.6667
$ws.Range("L113").Value = 6374.6665
$ws.Range("M113").Value = -921.6667000000002
$ws.Range("N113").Value = -10714.6665

$ws.Range("H122").Value = 94851.63
$ws.Range("I122").Value = 1604
$ws.Range("J122").Value = 129819.5
$ws.Range("K122").Value = 4812
$ws.Range("L122").Value = 389458.5
$ws.Range("M122").Value = -2362
$ws.Range("N122").Value = -394358.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 11198.2
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 11198.2
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 33594.60000000001
$ws.Range("N39").Value = -34182.60000000001

$ws.Range("H98").Value = 1025.9231
$ws.Range("I98").Value = 893.1667
$ws.Range("J98").Value = 1139.7142
$ws.Range("K98").Value = 2679.5001
$ws.Range("L98").Value = 3419.1426
$ws.Range("M98").Value = -1181.5001
$ws.Range("N98").Value = -6415.142599999999

$ws.Range("H132").Value = 8497.1
$ws.Range("I132").Value = 3290.7058
$ws.Range("J132").Value = 15305.462
$ws.Range("K132").Value = 29616.3522
$ws.Range("L132").Value = 137749.158
$ws.Range("M132").Value = -27086.3522
$ws.Range("N132").Value = -142809.158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 133.53847
$ws.Range("I2").Value = 136.33333
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 136.33333
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -23.33332999999999
$ws.Range("N2").Value = -326

$ws.Range("H122").Value = 1577428
$ws.Range("I122").Value = 2014318.2
$ws.Range("J122").Value = 4622.9
$ws.Range("K122").Value = 6042954.6
$ws.Range("L122").Value = 13868.7
$ws.Range("M122").Value = -6040504.6
$ws.Range("N122").Value = -18768.7

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 7475.6
$ws.Range("I132").Value = 2708.7144
$ws.Range("J132").Value = 11646.625
$ws.Range("K132").Value = 8126.1432
$ws.Range("L132").Value = 34939.875
$ws.Range("M132").Value = -5596.1432
$ws.Range("N132").Value = -39999.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2703.182
$ws.Range("I22").Value = 1333.3334
$ws.Range("J22").Value = 3216.875
$ws.Range("K22").Value = 1333.3334
$ws.Range("L22").Value = 3216.875
$ws.Range("M22").Value = -1038.3334
$ws.Range("N22").Value = -3806.875

$ws.Range("H27").Value = 2703.182
$ws.Range("I27").Value = 1333.3334
$ws.Range("J27").Value = 3216.875
$ws.Range("K27").Value = 1333.3334
$ws.Range("L27").Value = 3216.875
$ws.Range("M27").Value = -1226.3334
$ws.Range("N27").Value = -3430.875

$ws.Range("H40").Value = 5295.9165
$ws.Range("I40").Value = 3058.1667
$ws.Range("J40").Value = 7533.6665
$ws.Range("K40").Value = 3058.1667
$ws.Range("L40").Value = 7533.6665
$ws.Range("M40").Value = -2922.1667
$ws.Range("N40").Value = -7805.6665

$ws.Range("H61").Value = 2634222.5
$ws.Range("I61").Value = 4000983.5
$ws.Range("J61").Value = 5836.231
$ws.Range("K61").Value = 4000983.5
$ws.Range("L61").Value = 5836.231
$ws.Range("M61").Value = -4000781.5
$ws.Range("N61").Value = -6240.231

$ws.Range("H93").Value = 7108.4165
$ws.Range("I93").Value = 6530.1
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 6530.1
$ws.Range("L93").Value = 10000
$ws.Range("M93").Value = -5282.1
$ws.Range("N93").Value = -12496

$ws.Range("H103").Value = 42595.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 42595.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 42595.5
$ws.Range("N103").Value = -44939.5

$ws.Range("H113").Value = 2634222.5
$ws.Range("I113").Value = 4000983.5
$ws.Range("J113").Value = 5836.231
$ws.Range("K113").Value = 4000983.5
$ws.Range("L113").Value = 5836.231
$ws.Range("M113").Value = -3998813.5
$ws.Range("N113").Value = -10176.231

$ws.Range("H122").Value = 4684.486
$ws.Range("I122").Value = 4026.32
$ws.Range("J122").Value = 6329.9
$ws.Range("K122").Value = 12078.96
$ws.Range("L122").Value = 18989.7
$ws.Range("M122").Value = -9628.960000000001
$ws.Range("N122").Value = -23889.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 1999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 1999
$ws.Range("N62").Value = -3247

$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 9995
$ws.Range("N65").Value = -16235

$ws.Range("H93").Value = 53942
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 53942
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 53942
$ws.Range("N93").Value = -58934

$ws.Range("H113").Value = 1273.7693
$ws.Range("I113").Value = 718.1875
$ws.Range("J113").Value = 2162.7
$ws.Range("K113").Value = 2154.5625
$ws.Range("L113").Value = 6488.099999999999
$ws.Range("M113").Value = 15.4375
$ws.Range("N113").Value = -10828.1

$ws.Range("H122").Value = 106275.35
$ws.Range("I122").Value = 189023.27
$ws.Range("J122").Value = 5139
$ws.Range("K122").Value = 567069.8099999999
$ws.Range("L122").Value = 15417
$ws.Range("M122").Value = -564619.8099999999
$ws.Range("N122").Value = -20317

$ws.Range("H126").Value = 6999.8
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6999.8
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 20999.4
$ws.Range("N126").Value = -25939.4
$ws.Range("M126").ClearContents()
